# edit.ps1 -- apply the commit's edit:
#   1. Split the "...around 2 milliseconds..." run into "...around 2 m" | "illiseconds..."
#      and plant the (collapsed) "_GoBack" bookmark at that split point, mirroring the
#      cursor position Word leaves behind after the last keystroke of an editing session.
#   2. Merge the " along with ... would have " / "a positive " runs (which used to be
#      split apart by the old "_GoBack" bookmark) back into a single run, and drop the
#      old bookmark, since "_GoBack" only ever tracks the single most-recent edit spot.
#
# Overall visible text is unchanged; only run/bookmark bookkeeping moves.

$d = $word.ActiveDocument

# ===========================================================================
# Part 2 first: merge "... would have " + "a positive " into one run, and
# remove the old "_GoBack" bookmark that used to sit between them.
# ===========================================================================

# Locate the (unique) run boundaries around the edit using stable anchor text,
# scoping each successive Find to start after the previous match so we land on
# the correct occurrence (some of these words/phrases recur earlier in the doc).
$rngSchlick = $d.Content
$findSchlick = $rngSchlick.Find
$findSchlick.Execute("Schlick", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rngApprox = $d.Range($rngSchlick.End, $d.Content.End)
$findApprox = $rngApprox.Find
$findApprox.Execute(" approximation", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startOfHaveRun = $rngApprox.End

$rngEffect = $d.Range($startOfHaveRun, $d.Content.End)
$findEffect = $rngEffect.Find
$findEffect.Execute("effect", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startOfEffectRun = $rngEffect.Start

# Guard both sides of the region we are about to touch with temporary bookmarks
# so that the runtime's run-coalescing pass cannot cascade into the neighbouring
# "approximation" / "effect" runs (bookmarks act as a hard stop for that pass).
$leftGuard = $d.Range($startOfHaveRun, $startOfHaveRun)
$d.Bookmarks.Add("ZZTMP_LEFT", $leftGuard)

$rightGuard = $d.Range($startOfEffectRun, $startOfEffectRun)
$d.Bookmarks.Add("ZZTMP_RIGHT", $rightGuard)

# The old "_GoBack" bookmark sits between the two runs we are merging -- drop it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Merge the two runs by replacing the whole span with identical text. A straight
# no-op assignment is elided by the runtime, so first swap in a placeholder and
# then swap the real text back in -- that forces the run structure to rebuild as
# a single run.
$mergeRange = $d.Content
$mergeFind = $mergeRange.Find
$mergeFind.Execute(" along with this finding a way to sample less from textures would have a positive ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$finalText = $mergeRange.Text
$mergeRange.Text = "ZZPLACEHOLDERZZ"

$placeholderRange = $d.Content
$placeholderFind = $placeholderRange.Find
$placeholderFind.Execute("ZZPLACEHOLDERZZ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholderRange.Text = $finalText

# Drop the temporary guard bookmarks -- pure bookmark deletion does not touch
# run text, so it will not re-trigger the coalescing pass.
$d.Bookmarks.Item("ZZTMP_LEFT").Delete()
$d.Bookmarks.Item("ZZTMP_RIGHT").Delete()

# ===========================================================================
# Part 1: split "...around 2 milliseconds per spot light." into
# "...around 2 m" | "illiseconds per spot light." and plant "_GoBack" there.
# ===========================================================================

$splitRange = $d.Content
$splitFind = $splitRange.Find
$splitFind.Execute("costing around 2 m", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitRange)
